# Update TPM-derived NATMI metrics (ligand/receptor expression + edge weights)
# per the new TPM recomputation (see commit message: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 1.006697333333333
$ws.Range("H2").Value2 = 3.020092
$ws.Range("I2").Value2 = 0.0001985651645046208
$ws.Range("J2").Value2 = 0.0001985651645046208
$ws.Range("M2").Value2 = 2.680851666666667
$ws.Range("N2").Value2 = 8.042555
$ws.Range("O2").Value2 = 0.1074910720871699
$ws.Range("P2").Value2 = 0.1074910720871699
$ws.Range("Q2").Value2 = 2.698806223895556
$ws.Range("R2").Value2 = 24.28925601506
$ws.Range("S2").Value2 = 0.00002134398241176695
$ws.Range("T2").Value2 = 0.00002134398241176694

# Row 3
$ws.Range("G3").Value2 = 1.006697333333333
$ws.Range("H3").Value2 = 3.020092
$ws.Range("I3").Value2 = 0.0001985651645046208
$ws.Range("J3").Value2 = 0.0001985651645046208
$ws.Range("M3").Value2 = 18.51427066666667
$ws.Range("O3").Value2 = 0.7423457357290222
$ws.Range("P3").Value2 = 0.7423457357290222
$ws.Range("Q3").Value2 = 18.63826690874489
$ws.Range("R3").Value2 = 167.744402178704
$ws.Range("S3").Value2 = 0.000147404003134337
$ws.Range("T3").Value2 = 0.000147404003134337

# Row 4
$ws.Range("G4").Value2 = 1.006697333333333
$ws.Range("H4").Value2 = 3.020092
$ws.Range("I4").Value2 = 0.0001985651645046208
$ws.Range("J4").Value2 = 0.0001985651645046208
$ws.Range("M4").Value2 = 3.745104
$ws.Range("N4").Value2 = 11.235312
$ws.Range("O4").Value2 = 0.1501631921838079
$ws.Range("P4").Value2 = 0.1501631921838079
$ws.Range("Q4").Value2 = 3.770186209856
$ws.Range("R4").Value2 = 33.931675888704
$ws.Range("S4").Value2 = 0.0000298171789585168
$ws.Range("T4").Value2 = 0.0000298171789585168

# Row 5
$ws.Range("I5").Value2 = 0.9806494927176636
$ws.Range("J5").Value2 = 0.9806494927176637
$ws.Range("M5").Value2 = 2.680851666666667
$ws.Range("N5").Value2 = 8.042555
$ws.Range("O5").Value2 = 0.1074910720871699
$ws.Range("P5").Value2 = 0.1074910720871699
$ws.Range("Q5").Value2 = 13328.53605519946
$ws.Range("R5").Value2 = 119956.8244967951
$ws.Range("S5").Value2 = 0.105411065313961
$ws.Range("T5").Value2 = 0.105411065313961

# Row 6
$ws.Range("I6").Value2 = 0.9806494927176636
$ws.Range("J6").Value2 = 0.9806494927176637
$ws.Range("M6").Value2 = 18.51427066666667
$ws.Range("O6").Value2 = 0.7423457357290222
$ws.Range("P6").Value2 = 0.7423457357290222
$ws.Range("Q6").Value2 = 92048.40655104816
$ws.Range("R6").Value2 = 828435.6589594334
$ws.Range("S6").Value2 = 0.7279809691637863
$ws.Range("T6").Value2 = 0.7279809691637864

# Row 7
$ws.Range("I7").Value2 = 0.9806494927176636
$ws.Range("J7").Value2 = 0.9806494927176637
$ws.Range("M7").Value2 = 3.745104
$ws.Range("N7").Value2 = 11.235312
$ws.Range("O7").Value2 = 0.1501631921838079
$ws.Range("P7").Value2 = 0.1501631921838079
$ws.Range("Q7").Value2 = 18619.73727048371
$ws.Range("R7").Value2 = 167577.6354343534
$ws.Range("S7").Value2 = 0.1472574582399162
$ws.Range("T7").Value2 = 0.1472574582399162

# Row 8
$ws.Range("G8").Value2 = 97.097641
$ws.Range("H8").Value2 = 291.292923
$ws.Range("I8").Value2 = 0.01915194211783179
$ws.Range("J8").Value2 = 0.01915194211783179
$ws.Range("M8").Value2 = 2.680851666666667
$ws.Range("N8").Value2 = 8.042555
$ws.Range("O8").Value2 = 0.1074910720871699
$ws.Range("P8").Value2 = 0.1074910720871699
$ws.Range("Q8").Value2 = 260.3043727042516
$ws.Range("R8").Value2 = 2342.739354338265
$ws.Range("S8").Value2 = 0.002058662790797162
$ws.Range("T8").Value2 = 0.002058662790797162

# Row 9
$ws.Range("G9").Value2 = 97.097641
$ws.Range("H9").Value2 = 291.292923
$ws.Range("I9").Value2 = 0.01915194211783179
$ws.Range("J9").Value2 = 0.01915194211783179
$ws.Range("M9").Value2 = 18.51427066666667
$ws.Range("O9").Value2 = 0.7423457357290222
$ws.Range("P9").Value2 = 0.7423457357290222
$ws.Range("Q9").Value2 = 1797.692006568831
$ws.Range("R9").Value2 = 16179.22805911948
$ws.Range("S9").Value2 = 0.01421736256210149
$ws.Range("T9").Value2 = 0.01421736256210149

# Row 10
$ws.Range("G10").Value2 = 97.097641
$ws.Range("H10").Value2 = 291.292923
$ws.Range("I10").Value2 = 0.01915194211783179
$ws.Range("J10").Value2 = 0.01915194211783179
$ws.Range("M10").Value2 = 3.745104
$ws.Range("N10").Value2 = 11.235312
$ws.Range("O10").Value2 = 0.1501631921838079
$ws.Range("P10").Value2 = 0.1501631921838079
$ws.Range("Q10").Value2 = 363.640763699664
$ws.Range("R10").Value2 = 3272.766873296976
$ws.Range("S10").Value2 = 0.002875916764933139
$ws.Range("T10").Value2 = 0.002875916764933139
